$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Row 2 (TC1): replace the ERR-5005 details with the ERR-1001 details
$ws.Range("B2").Value = "ERR-1001"
$ws.Range("C2").Value = "connection timeout"
$ws.Range("D2").Value = "WARN - ERR-1001: connection timeout in module Fixflex. User: nicky.stracke, SessionID: b92e0e36-d08c-4864-b072-bcf6b0a5868f"

# Row 3 (TC2): replace the ERR-5005 details with the new ERR-3003 details
$ws.Range("B3").Value = "ERR-3003"
$ws.Range("C3").Value = "Invalid input parameter"
$ws.Range("D3").Value = "WARN - ERR-3003: Invalid input parameter in module Tempsoft. User: kyle.windler, SessionID: 5f98e075-91e7-4882-bec8-79b4bb6920ea"
